$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.534.96'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '1.658.58'
$ws.Range('E3').Value = '  -3.07%  '
$ws.Range('E4').Value = '  +0.97%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.68'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('E6').Value = '  -1.11%  '
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.23'
$ws.Range('E8').Value = '  -2.40%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.259'
$ws.Range('E9').Value = '  -1.90%  '
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('E11').Value = '  -1.94%  '
$ws.Range('D12').Value = '1.892.65'
$ws.Range('E12').Value = '  -2.91%  '
$ws.Range('D13').Value = '1.656.71'
$ws.Range('E13').Value = '  -3.05%  '
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('E15').Value = '  -3.50%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.72'
$ws.Range('E16').Value = '  -2.50%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '245.46'
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('D18').Value = '27.509.73'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').Value = '0.0₃0729'
$ws.Range('E19').Value = '  -2.35%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.48'
$ws.Range('E20').Value = '  -6.16%  '
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.47'
$ws.Range('E22').Value = '  -2.65%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.30'
$ws.Range('E23').Value = '  -3.79%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  -3.13%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.98'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('E26').Value = '  -5.24%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.20'
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('E28').Value = '  +1.13%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.111'
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('E30').Value = '  +5.24%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0499'
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('E32').Value = '  -2.02%  '
$ws.Range('D33').Value = '1.435.63'
$ws.Range('E33').Value = '  -7.27%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.14'
$ws.Range('E34').Value = '  -4.14%  '
$ws.Range('E35').Value = '  -8.43%  '
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.930'
$ws.Range('E37').Value = '  -4.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.580'
$ws.Range('E38').Value = '  -5.69%  '
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.04'
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '69.12'
$ws.Range('E41').Value = '  -2.81%  '
$ws.Range('E42').Value = '  +1.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.39'
$ws.Range('E43').Value = '  -7.96%  '
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.801.44'
$ws.Range('E45').Value = '  -2.15%  '
$ws.Range('B46').Value = 'MXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.20'
$ws.Range('E46').Value = '  -6.16%  '
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '88.79'
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('D49').Value = '0.0₆0108'
$ws.Range('E49').Value = '  +4.51%  '
$ws.Range('E50').Value = '  -4.92%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.81'
$ws.Range('E51').Value = '  -6.89%  '
